# Auto-generated edit script for cs-en-us-pbms.xlsx weekly update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared rich-text strings) ---
# "Volume 30   Number  25" -> "...  26"  (A8)
$ws.Range("A8").Characters(21, 2).Text = "26"

# "Report Covering the Week  6/19/2023  Through  6/25/2023" -> new dates (C9)
# Replace the later substring first so offsets of the earlier one stay valid.
$ws.Range("C9").Characters(47, 9).Text = "7/2/2023"
$ws.Range("C9").Characters(27, 9).Text = "6/26/2023"

# --- Crime-statistics table updates (rows 14-30) ---
$ws.Range("D14").Copy($ws.Range("C14"))
$ws.Range("N14").Value = -79.545454545454
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 13
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 8.333333333333
$ws.Range("I15").Value = 56
$ws.Range("J15").Value = 87
$ws.Range("K15").Value = -35.632183908046
$ws.Range("L15").Value = -25.333333333333
$ws.Range("M15").Value = 27.272727272727
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = 38
$ws.Range("D16").Value = 59
$ws.Range("E16").Value = -35.593220338983
$ws.Range("F16").Value = 151
$ws.Range("G16").Value = 184
$ws.Range("H16").Value = -17.934782608695
$ws.Range("I16").Value = 877
$ws.Range("J16").Value = 1032
$ws.Range("K16").Value = -15.019379844961
$ws.Range("L16").Value = 31.091180866965
$ws.Range("M16").Value = 52.256944444444
$ws.Range("N16").Value = -82.847643262272
$ws.Range("C17").Value = 48
$ws.Range("D17").Value = 45
$ws.Range("E17").Value = 6.666666666666
$ws.Range("F17").Value = 185
$ws.Range("G17").Value = 183
$ws.Range("H17").Value = 1.092896174863
$ws.Range("I17").Value = 1070
$ws.Range("J17").Value = 983
$ws.Range("K17").Value = 8.850457782299
$ws.Range("L17").Value = 17.971334068357
$ws.Range("M17").Value = 68.769716088328
$ws.Range("N17").Value = -36.385255648038
$ws.Range("C18").Value = 26
$ws.Range("D18").Value = 56
$ws.Range("E18").Value = -53.571428571428
$ws.Range("F18").Value = 132
$ws.Range("G18").Value = 265
$ws.Range("H18").Value = -50.188679245283
$ws.Range("I18").Value = 1060
$ws.Range("J18").Value = 1517
$ws.Range("K18").Value = -30.125247198417
$ws.Range("L18").Value = 14.594594594594
$ws.Range("M18").Value = 15.217391304347
$ws.Range("N18").Value = -81.932844724731
$ws.Range("C19").Value = 214
$ws.Range("D19").Value = 255
$ws.Range("E19").Value = -16.078431372549
$ws.Range("F19").Value = 952
$ws.Range("G19").Value = 982
$ws.Range("H19").Value = -3.054989816700
$ws.Range("I19").Value = 5715
$ws.Range("J19").Value = 5483
$ws.Range("K19").Value = 4.231260258982
$ws.Range("L19").Value = 77.594779366065
$ws.Range("M19").Value = 12.014896119169
$ws.Range("N19").Value = -63.947766843300
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 24
$ws.Range("E20").Value = -41.666666666666
$ws.Range("F20").Value = 49
$ws.Range("G20").Value = 92
$ws.Range("H20").Value = -46.739130434782
$ws.Range("I20").Value = 275
$ws.Range("J20").Value = 332
$ws.Range("K20").Value = -17.168674698795
$ws.Range("L20").Value = 7.003891050583
$ws.Range("M20").Value = 51.933701657458
$ws.Range("N20").Value = -90.959894806048
$ws.Range("C21").Value = 343
$ws.Range("D21").Value = 442
$ws.Range("E21").Value = -22.398190045248
$ws.Range("F21").Value = 1486
$ws.Range("G21").Value = 1718
$ws.Range("H21").Value = -13.504074505238
$ws.Range("I21").Value = 9062
$ws.Range("J21").Value = 9452
$ws.Range("K21").Value = -4.126110876005
$ws.Range("L21").Value = 49.439313984168
$ws.Range("M21").Value = 21.409431939978
$ws.Range("N21").Value = -71.424066599394
$ws.Range("C22").Value = 19
$ws.Range("D22").Value = 11
$ws.Range("E22").Value = 72.727272727272
$ws.Range("F22").Value = 57
$ws.Range("G22").Value = 44
$ws.Range("H22").Value = 29.545454545454
$ws.Range("I22").Value = 327
$ws.Range("J22").Value = 348
$ws.Range("K22").Value = -6.034482758620
$ws.Range("L22").Value = 46.636771300448
$ws.Range("M22").Value = 26.744186046511
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 80
$ws.Range("G23").Value = 29
$ws.Range("H23").Value = -13.793103448275
$ws.Range("I23").Value = 190
$ws.Range("J23").Value = 227
$ws.Range("K23").Value = -16.299559471365
$ws.Range("L23").Value = -25.490196078431
$ws.Range("M23").Value = -0.523560209424
$ws.Range("C24").Value = 446
$ws.Range("D24").Value = 465
$ws.Range("E24").Value = -4.086021505376
$ws.Range("F24").Value = 1808
$ws.Range("G24").Value = 1877
$ws.Range("H24").Value = -3.676078849227
$ws.Range("I24").Value = 10034
$ws.Range("J24").Value = 10590
$ws.Range("K24").Value = -5.250236071765
$ws.Range("L24").Value = 56.317183361894
$ws.Range("M24").Value = 20.196454240536
$ws.Range("C25").Value = 118
$ws.Range("D25").Value = 101
$ws.Range("E25").Value = 16.831683168316
$ws.Range("F25").Value = 431
$ws.Range("G25").Value = 396
$ws.Range("H25").Value = 8.838383838383
$ws.Range("I25").Value = 2477
$ws.Range("J25").Value = 2278
$ws.Range("K25").Value = 8.735733099209
$ws.Range("L25").Value = 41.542857142857
$ws.Range("M25").Value = 36.323610346725
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 14
$ws.Range("H26").Value = -17.647058823529
$ws.Range("I26").Value = 102
$ws.Range("J26").Value = 138
$ws.Range("K26").Value = -26.086956521739
$ws.Range("L26").Value = -12.068965517241
$ws.Range("C27").Value = 25
$ws.Range("D27").Value = 25
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 88
$ws.Range("H27").Value = -7.368421052631
$ws.Range("I27").Value = 464
$ws.Range("J27").Value = 470
$ws.Range("K27").Value = -1.276595744680
$ws.Range("L27").Value = 29.247910863509
$ws.Range("C28").Value = 2
$ws.Range("F14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 3
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -33.333333333333
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 21
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = -30
$ws.Range("L28").Value = -19.230769230769
$ws.Range("M28").Value = 61.538461538461
$ws.Range("N28").Value = -70
$ws.Range("F14").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 2
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -50
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 18
$ws.Range("J29").Value = 25
$ws.Range("K29").Value = -28
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -70.967741935483
$ws.Range("D14").Copy($ws.Range("C30"))
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 3
$ws.Range("H30").Value = -80
$ws.Range("I30").Value = 44
$ws.Range("J30").Value = 92
$ws.Range("K30").Value = -52.173913043478
$ws.Range("L30").Value = -48.235294117647
